# Update keyboard type names in the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the keyboard type values in column C (rows 2-28):
#   "Rubber Dome" -> "Dome-Switch"
#   "Membrane"    -> "Scissor-Switch"
#   "Mechanical"  stays "Mechanical"
$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Text -eq "Rubber Dome") {
        $cell.Value = "Dome-Switch"
    } elseif ($cell.Text -eq "Membrane") {
        $cell.Value = "Scissor-Switch"
    }
}

# Reset the view so the top-left visible cell is back to A1.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
